# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404   (left-hand "before" columns, A:J)
#   *_new -> *_FV2410   (right-hand "after" columns, L:U)
# then wrap the sheet's data range in a table and freeze the header row,
# matching the regenerated merged AHB export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" / "_new" header suffixes -------------------------
# Columns A-J carry the "_old" suffixed headers (column K is the "diff"
# column and stays untouched), columns L-U carry the "_new" suffixed ones.
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value() -replace "_old$", "_FV2404")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value() -replace "_new$", "_FV2410")
}

# --- 2. Turn the used range into a native Excel table -----------------------
$dataRange = $ws.Range("A1:U61")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
